# Update imputed values in the RandomForest result data worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.989900000000006
$ws.Range("D18").Value = -8.576499999999999
$ws.Range("B21").Value = 8.883800000000001
$ws.Range("B23").Value = 9.214499999999997
$ws.Range("C24").Value = -13.33879999999999
$ws.Range("B25").Value = 5.817200000000001
$ws.Range("C28").Value = -13.41009999999999
$ws.Range("C36").Value = -11.75010000000001
$ws.Range("C45").Value = -14.01329999999999
$ws.Range("C48").Value = -11.2821
$ws.Range("C49").Value = -13.8087
$ws.Range("D51").Value = -8.3522
$ws.Range("C52").Value = -10.79059999999999
$ws.Range("B53").Value = 6.414999999999996
$ws.Range("C53").Value = -10.67410000000001
$ws.Range("C54").Value = -13.50390000000001
$ws.Range("D55").Value = -8.6159
$ws.Range("B57").Value = 4.828699999999998
$ws.Range("B59").Value = 5.041599999999999
$ws.Range("D64").Value = -7.840199999999989
$ws.Range("B69").Value = 5.732199999999994
$ws.Range("C70").Value = -11.5357
$ws.Range("B79").Value = 9.314700000000004
$ws.Range("D80").Value = -8.100099999999999
$ws.Range("B83").Value = 5.182099999999999
$ws.Range("C86").Value = -13.84169999999999
$ws.Range("C87").Value = -13.351
$ws.Range("D92").Value = -6.989900000000005
$ws.Range("B93").Value = 5.776899999999999
$ws.Range("D94").Value = -6.231800000000003
$ws.Range("D96").Value = -8.598999999999998
$ws.Range("C101").Value = -12.89970000000001
